# VAT of PS added in Bill
# Applies the changes described by the commit to Sheet1:
#  - C6 description updated to "PS including VAT"
#  - F6:F14 style cleanup (drop the redundant applyFill style, matches style used by E/G columns)
#  - Row 16 (G/H/I) rewritten as a (shared) formula 1.5+0.1+5
#  - Two new rows (18 "VAT Amount of PS" / 19 "Deduction% in VAT Amount") added
#  - Selection moved to B20

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Clean up F6:F14 cell style (was using the redundant "applyFill" border-only
#        style; align it with the plain bordered style already used by column E/G) ---
$ws.Range("E6").Copy()
$ws.Range("F6:F14").PasteSpecial(-4122)

# Re-assert the text values for F6:F14 (paste-special of formats only shouldn't have
# touched them, but make sure they still read correctly / shared strings stay intact)
$ws.Range("F6").Value2 = "PS"
$ws.Range("F7").Value2 = "Subtotal"
$ws.Range("F8").Value2 = "VAT %"
$ws.Range("F9").Value2 = "VAT Amount"
$ws.Range("F10").Value2 = "Contingency %"
$ws.Range("F11").Value2 = "Physical Contingency %"
$ws.Range("F12").Value2 = "Price Contingency %"
$ws.Range("F13").Value2 = "Total (A+B+D)"
$ws.Range("F14").Value2 = "GrandTotal incl. contingencies"

# --- 2. New row 18 content/format (pattern copied from row 6: D/H highlighted input
#        cells, everything else plain bordered) ---
$ws.Range("B6:I6").Copy()
$ws.Range("B18:I18").PasteSpecial(-4122)

# --- 3. New row 19 content/format (pattern copied from row 17: everything plain
#        bordered, no highlighted input cells) ---
$ws.Range("B17:I17").Copy()
$ws.Range("B19:I19").PasteSpecial(-4122)

# --- 4. Fill in the new row values, in the same order the strings were authored so
#        the shared-string table ends up matching (M/VAT Amount of PS, then the C6
#        rename, then Deduction%/N) ---
$ws.Range("A18").Value2 = 12
$ws.Range("B18").Value2 = "M"
$ws.Range("C18").Value2 = "VAT Amount of PS"
$ws.Range("D18").Value2 = ""
$ws.Range("E18").Value2 = ""
$ws.Range("F18").Value2 = "VAT Amount of PS"
$ws.Range("G18").Value2 = 0
$ws.Range("H18").Value2 = ""
$ws.Range("I18").Value2 = ""

# --- 5. C6 description rename ---
$ws.Range("C6").Value2 = "PS including VAT"

# --- 6. Row 19 values ---
$ws.Range("A19").Value2 = 13
$ws.Range("C19").Value2 = "Deduction% in VAT Amount"
$ws.Range("F19").Value2 = "Deduction% in VAT Amount"
$ws.Range("B19").Value2 = "N"
$ws.Range("D19").Value2 = 0
$ws.Range("E19").Value2 = 0
$ws.Range("G19").Value2 = 30
$ws.Range("H19").Value2 = 30
$ws.Range("I19").Value2 = 30

# --- 7. Row 16 G/H/I rewritten as formulas (H/I become a shared formula group) ---
$ws.Range("G16").Formula = "=1.5+0.1+5"
$ws.Range("H16:I16").Formula = "=1.5+0.1+5"

# --- 8. Move the active selection to B20 ---
[void]$ws.Range("B20").Select()
